# "cleaned data with filters and got status at good speed!"
#
# Fills in the missing "Status" (column F) value for the 30 rows that were
# left blank, matching the existing "Active" / "Open" values already used
# throughout column F (shared-string table entries are reused verbatim,
# including their exact whitespace, so Excel doesn't mint new duplicate
# strings). Each of those rows also grows to row height 60 (matching every
# other wrapped Status cell), column D is widened to fit its contents, and
# the sheet's scroll/selection state is refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact text of the two shared strings already used for "Active" / "Open"
# (indices 579 / 580 in sharedStrings.xml) so Excel reuses them instead of
# creating new duplicate entries.
$active = "`n                    Active`n                "
$openStatus = "`n                  Open`n              "

# row number -> status ("Active" or "Open")
$statusByRow = [ordered]@{
    16  = $openStatus
    35  = $active
    42  = $active
    57  = $openStatus
    60  = $openStatus
    68  = $active
    72  = $active
    75  = $openStatus
    83  = $openStatus
    91  = $openStatus
    95  = $active
    106 = $openStatus
    128 = $active
    130 = $openStatus
    160 = $openStatus
    223 = $active
    226 = $openStatus
    243 = $openStatus
    247 = $active
    255 = $openStatus
    258 = $active
    259 = $active
    260 = $openStatus
    261 = $active
    262 = $openStatus
    263 = $openStatus
    264 = $active
    265 = $openStatus
    266 = $active
    267 = $active
}

foreach ($row in $statusByRow.Keys) {
    $cell = $ws.Range("F$row")
    $cell.WrapText = $true
    $cell.Value = $statusByRow[$row]
    $ws.Rows.Item($row).RowHeight = 60
}

# Column D ("BR number") now has values long enough to need a wider, best-fit
# column (roughly matching the widest entry, e.g. "381502.00 mln GBP").
$ws.Columns.Item(4).ColumnWidth = 16.8

# Refresh the view: scroll back to the top and select the last-edited cell.
$ws.Range("E223").Select()
